# Update the "want to go" counts (column F) for a handful of events on the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet, which
# mirrors the same events. The site was regenerated, bumping a few counters.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 11699
$wsExhibit.Range("F7").Value = 11654
$wsExhibit.Range("F14").Value = 3517

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 11699
$wsAll.Range("F9").Value = 11654
$wsAll.Range("F17").Value = 3517
